# Auto-update draw results: append the 2025-10-30 Pick 3 row to the Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$row = 44

# New row values (all stored as text, matching the existing table's columns).
# A leading apostrophe forces Excel to keep number/date-looking text ("2025-10-30",
# "251030", the ISO timestamp) as literal text instead of auto-converting it to a
# date serial / numeric value - same as every other row already in the sheet.
$ws.Range("A" + $row).Value = "'2025-10-30"
$ws.Range("B" + $row).Value = "Pick 3"
$ws.Range("C" + $row).Value = "'251030"
$ws.Range("D" + $row).Value = "7-2-3"
$ws.Range("E" + $row).Value = "'2025-10-30T21:40:11.797+04:00"
